$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column cells being updated keep their original text
# representation (Excel would otherwise auto-convert plain numeric-looking
# strings like "212.77" into actual numbers).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D12", "D13", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D25", "D26", "D30", "D33", "D37", "D38", "D40", "D41", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.662.23"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.637.94"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "212.77"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "23.07"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "0.0610"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.870.68"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.636.85"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "64.59"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "27.659.00"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "230.00"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "7.73"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "0.0₃0722"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "10.03"
$ws.Range("E23").Value = "  +4.05%  "
$ws.Range("E24").Value = "  -2.91%  "
$ws.Range("D25").Value = "149.91"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").Value = "6.92"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "1.449.95"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "0.565"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "0.876"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "0.904"
$ws.Range("E40").Value = "  +10.54%  "
$ws.Range("D41").Value = "69.84"
$ws.Range("E41").Value = "  +8.11%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "1.780.53"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "1.71"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "86.17"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "0.0₆0107"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "0.0987"
$ws.Range("E51").Value = "  -0.36%  "
